# Auto commit at 2025-09-09  8:15:51.90
#
# Refresh the raw metric figures on the "Metrics" sheet (B2:B13). The
# "today" sheet pulls these via =Metrics!Bn formulas (and downstream
# =B.. / =E..+B.. running totals), so it recalculates automatically once
# the source cells change - no need to touch it directly except for its
# selection.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value  = 127298.02999999998
$metrics.Range("B3").Value  = 103357.43000000002
$metrics.Range("B4").Value  = 40384.850000000006
$metrics.Range("B5").Value  = 4975
$metrics.Range("B6").Value  = 4046548.9099999997
$metrics.Range("B7").Value  = 3430884.9099999992
$metrics.Range("B8").Value  = 1169750.53
$metrics.Range("B9").Value  = 156135
$metrics.Range("B10").Value = 32511872.710999828
$metrics.Range("B11").Value = 19460754.980000004
$metrics.Range("B12").Value = 11451459.420000002
$metrics.Range("B13").Value = 1253762

# Selection on "Metrics" moved from D10 to D43.
$metrics.Activate()
$metrics.Range("D43").Select()

# Selection on "today" (the tab that stays active/selected) moved from
# E8 to G21.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("G21").Select()
